$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Points for grading" (column E) scores that were left blank,
# matching the full marks already recorded in column D (Total Points),
# for the "Generic" and "Customer Class" rubric sections.
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Move the selection/active cell to E15, matching where the grader left off.
[void]$ws.Range("E15").Select()
